$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Number of used rows (header + data) before we touch anything.
$lastRow = $ws.UsedRange.Rows.Count

# Insert a new column ahead of the current "Polarity" column (B), which
# shifts Polarity -> C and Review -> D. This becomes the "Unnamed: 0.1"
# column, a duplicate of the original "Unnamed: 0" index column.
$ws.Columns("B").Insert()

# Give the new header cell the same (bold/centered) look as its neighbours,
# then set its text.
$ws.Range("A1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$ws.Range("B1").Value = "Unnamed: 0.1"

# Fill the new column with the same values as column A ("Unnamed: 0").
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells($r, 2).Value = $ws.Cells($r, 1).Value2
}

# Lowercase the review text, which now lives in column D.
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells($r, 4)
    $cell.Value = $cell.Value2.ToLower()
}
